$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header cell F4 style (drop the duplicate style, reuse the one used by E4) ---
$ws.Range("E4").Copy()
$ws.Range("F4").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 6: correct the logged start/end time (12:00 PM - 6:00 PM instead of 1:00 PM - 5:00 PM) ---
$ws.Range("G6").Value = 0.5
$ws.Range("H6").Value = 0.75

# --- Row 7: add the new "week 3" log entry ---
$ws.Rows(7).RowHeight = 62.4

$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "week 3"
$ws.Range("C7").Value = "fri"

$ws.Range("D6").Copy()
$ws.Range("D7").PasteSpecial(-4122)   # xlPasteFormats (date format)
$ws.Range("D7").Value = 45863

$ws.Range("E7").Value = "offline"

$ws.Range("F6").Copy()
$ws.Range("F7").PasteSpecial(-4122)   # xlPasteFormats (wrap text)
$ws.Range("F7").Value = "observed the HMI system. Took insights . Brainstromed solution for inventory mangement. Developed basic logic for Pill distribution . "

$ws.Range("G7").Value = 0.375
$ws.Range("H7").Value = 0.66666666666666663

# --- Update selection to match the author's last cursor position ---
$ws.Range("E6").Select()

$excel.CalculateFull()
